$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new rows (31 and 32) following the existing pattern.
$ws.Cells.Item(31, 1).Value = 10001
$ws.Cells.Item(31, 2).Value = 110030
$ws.Cells.Item(31, 3).Value = "eng"
$ws.Cells.Item(31, 4).Value = $true
$ws.Cells.Item(31, 5).Value = "superadmin"
$ws.Cells.Item(31, 6).Value = "now()"

$ws.Cells.Item(32, 1).Value = 10001
$ws.Cells.Item(32, 2).Value = 110031
$ws.Cells.Item(32, 3).Value = "eng"
$ws.Cells.Item(32, 4).Value = $true
$ws.Cells.Item(32, 5).Value = "superadmin"
$ws.Cells.Item(32, 6).Value = "now()"

# Update the selection to match the post-edit view state (cursor on D32).
$ws.Range("D32").Select()
